# "multi OS Test Options and Test Properties"
# Adds two new test-data sheets (SignUp, FoodSearch) after the existing
# "Login" sheet, and updates one browser value on the Login sheet.

$wb    = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("Login")

# --- Login sheet: row 4 browser safari -> chrome ---
$login.Range("B4").Value = "chrome"

# --- New sheet: SignUp (placed right after Login) ---
$signup = $wb.Worksheets.Add($null, $login)
$signup.Name = "SignUp"

$signup.Range("A1").Value = "ID"
$signup.Range("B1").Value = "browser"
$signup.Range("C1").Value = "remote"
$signup.Range("A2").Value = "'9757078221"
$signup.Range("B2").Value = "firefox"
$signup.Range("C2").Value = "n"

# --- New sheet: FoodSearch (placed right after SignUp) ---
$foodsearch = $wb.Worksheets.Add($null, $signup)
$foodsearch.Name = "FoodSearch"

$foodsearch.Range("A1").Value = "location"
$foodsearch.Range("B1").Value = "browser"
$foodsearch.Range("A2").Value = "Thane"
$foodsearch.Range("B2").Value = "chrome"
$foodsearch.Range("A3").Value = "Vashi"
$foodsearch.Range("B3").Value = "firefox"

# Leave each sheet's selection where data entry ended, then return focus
# to the Login tab (matches the saved workbook's active-sheet state).
$null = $foodsearch.Range("B3").Select()
$null = $signup.Range("B2").Select()
$null = $login.Range("B4").Select()
$null = $login.Activate()
